# Auto-generated Excel COM-interop script
# Updates cached market-price / profit figures on each class sheet
# (currentAveragePrice*, LevePrice*, LeveProfit* columns), matching a
# scheduled market-data refresh. No formulas are involved -- all of
# these cells hold static cached numbers.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 213
$ws.Range("I9").Value = 230.83333
$ws.Range("J9").Value = 159.5
$ws.Range("K9").Value = 230.83333
$ws.Range("L9").Value = 159.5
$ws.Range("M9").Value = -61.83332999999999
$ws.Range("N9").Value = -497.5
$ws.Range("H43").Value = 7540.2
$ws.Range("I43").Value = 8300.25
$ws.Range("K43").Value = 8300.25
$ws.Range("M43").Value = -8231.25
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H86").Value = 3657.6667
$ws.Range("I86").Value = 3244
$ws.Range("J86").Value = 4485
$ws.Range("K86").Value = 3244
$ws.Range("L86").Value = 4485
$ws.Range("M86").Value = -2121
$ws.Range("N86").Value = -6731
$ws.Range("H89").Value = 3657.6667
$ws.Range("I89").Value = 3244
$ws.Range("J89").Value = 4485
$ws.Range("K89").Value = 16220
$ws.Range("L89").Value = 22425
$ws.Range("M89").Value = -10604
$ws.Range("N89").Value = -33657
$ws.Range("H101").Value = 565.1429000000001
$ws.Range("I101").Value = 576
$ws.Range("K101").Value = 1728
$ws.Range("M101").Value = -106
$ws.Range("H132").Value = 1101.1364
$ws.Range("I132").Value = 1150.7
$ws.Range("K132").Value = 3452.1
$ws.Range("M132").Value = -922.1000000000004
$ws.Range("H135").Value = 1059
$ws.Range("I135").Value = 951.8333
$ws.Range("K135").Value = 8566.4997
$ws.Range("M135").Value = -6031.4997
$ws.Range("H137").Value = 1769.5625
$ws.Range("J137").Value = 1794.3334
$ws.Range("L137").Value = 5383.0002
$ws.Range("N137").Value = -10483.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H45").Value = 27312.273
$ws.Range("I45").Value = 804.75
$ws.Range("K45").Value = 804.75
$ws.Range("M45").Value = -427.75
$ws.Range("H106").Value = 14999.5
$ws.Range("J106").Value = 14999.5
$ws.Range("L106").Value = 14999.5
$ws.Range("N106").Value = -17523.5
$ws.Range("H119").Value = 29999
$ws.Range("J119").Value = 29999
$ws.Range("L119").Value = 29999
$ws.Range("N119").Value = -39675
$ws.Range("H122").Value = 5708.278
$ws.Range("I122").Value = 5020.75
$ws.Range("J122").Value = 7083.3335
$ws.Range("K122").Value = 15062.25
$ws.Range("L122").Value = 21250.0005
$ws.Range("M122").Value = -12612.25
$ws.Range("N122").Value = -26150.0005
$ws.Range("H132").Value = 996.64703
$ws.Range("I132").Value = 996.64703
$ws.Range("K132").Value = 2989.94109
$ws.Range("M132").Value = -459.9410899999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1722.5454
$ws.Range("J134").Value = 1000
$ws.Range("L134").Value = 3000
$ws.Range("N134").Value = -8070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 138.375
$ws.Range("I7").Value = 77.333336
$ws.Range("J7").Value = 321.5
$ws.Range("K7").Value = 77.333336
$ws.Range("L7").Value = 321.5
$ws.Range("M7").Value = 35.666664
$ws.Range("N7").Value = -547.5
$ws.Range("H31").Value = 3783.9092
$ws.Range("I31").Value = 3270.3157
$ws.Range("K31").Value = 3270.3157
$ws.Range("M31").Value = -2975.3157
$ws.Range("H34").Value = 3783.9092
$ws.Range("I34").Value = 3270.3157
$ws.Range("K34").Value = 3270.3157
$ws.Range("M34").Value = -3068.3157
$ws.Range("H58").Value = 1932.4849
$ws.Range("J58").Value = 6493.2
$ws.Range("L58").Value = 6493.2
$ws.Range("N58").Value = -6899.2
$ws.Range("H86").Value = 6333.3335
$ws.Range("I86").Value = 4200
$ws.Range("J86").Value = 17000
$ws.Range("K86").Value = 4200
$ws.Range("L86").Value = 17000
$ws.Range("M86").Value = -3077
$ws.Range("N86").Value = -19246
$ws.Range("H89").Value = 6333.3335
$ws.Range("I89").Value = 4200
$ws.Range("J89").Value = 17000
$ws.Range("K89").Value = 21000
$ws.Range("L89").Value = 85000
$ws.Range("M89").Value = -15384
$ws.Range("N89").Value = -96232
$ws.Range("H122").Value = 2202.9048
$ws.Range("I122").Value = 2289.923
$ws.Range("J122").Value = 2061.5
$ws.Range("K122").Value = 6869.768999999999
$ws.Range("L122").Value = 6184.5
$ws.Range("M122").Value = -4419.768999999999
$ws.Range("N122").Value = -11084.5
$ws.Range("H134").Value = 3980.8333
$ws.Range("I134").Value = 3971.75
$ws.Range("K134").Value = 11915.25
$ws.Range("M134").Value = -9380.25
$ws.Range("H136").Value = 1932.4849
$ws.Range("J136").Value = 6493.2
$ws.Range("L136").Value = 19479.6
$ws.Range("N136").Value = -24579.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 175.85715
$ws.Range("I23").Value = 175.85715
$ws.Range("K23").Value = 527.5714499999999
$ws.Range("M23").Value = -292.5714499999999
$ws.Range("H62").Value = 3656.25
$ws.Range("J62").Value = 2964.2856
$ws.Range("L62").Value = 8892.856800000001
$ws.Range("N62").Value = -10264.8568
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H65").Value = 3656.25
$ws.Range("J65").Value = 2964.2856
$ws.Range("L65").Value = 26678.5704
$ws.Range("N65").Value = -33542.5704
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H88").Value = 4904.6665
$ws.Range("I88").Value = 4904.6665
$ws.Range("K88").Value = 14713.9995
$ws.Range("M88").Value = -14285.9995
$ws.Range("H91").Value = 4904.6665
$ws.Range("I91").Value = 4904.6665
$ws.Range("K91").Value = 14713.9995
$ws.Range("M91").Value = -13231.9995
$ws.Range("H113").Value = 3550.8
$ws.Range("J113").Value = 3213.75
$ws.Range("L113").Value = 9641.25
$ws.Range("N113").Value = -13981.25
$ws.Range("H131").Value = 1593.9584
$ws.Range("J131").Value = 1618.4783
$ws.Range("L131").Value = 4855.4349
$ws.Range("N131").Value = -14935.4349

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 62.5
$ws.Range("H97").Value = 993.3333
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H102").Value = 3004
$ws.Range("I102").Value = 2012
$ws.Range("K102").Value = 2012
$ws.Range("M102").Value = -390
$ws.Range("H122").Value = 95784.91
$ws.Range("I122").Value = 5133.2
$ws.Range("K122").Value = 15399.6
$ws.Range("M122").Value = -12949.6
$ws.Range("H126").Value = 5006.25
$ws.Range("I126").Value = 4005.5
$ws.Range("K126").Value = 12016.5
$ws.Range("M126").Value = -9546.5
$ws.Range("H132").Value = 3214.3845
$ws.Range("I132").Value = 2978.7
$ws.Range("K132").Value = 8936.099999999999
$ws.Range("M132").Value = -6406.099999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2362.2092
$ws.Range("I22").Value = 1689.0667
$ws.Range("J22").Value = 3915.6155
$ws.Range("K22").Value = 1689.0667
$ws.Range("L22").Value = 3915.6155
$ws.Range("M22").Value = -1394.0667
$ws.Range("N22").Value = -4505.6155
$ws.Range("H27").Value = 2362.2092
$ws.Range("I27").Value = 1689.0667
$ws.Range("J27").Value = 3915.6155
$ws.Range("K27").Value = 1689.0667
$ws.Range("L27").Value = 3915.6155
$ws.Range("M27").Value = -1582.0667
$ws.Range("N27").Value = -4129.6155
$ws.Range("H55").Value = 899.9
$ws.Range("I55").Value = 687.375
$ws.Range("J55").Value = 1750
$ws.Range("K55").Value = 687.375
$ws.Range("L55").Value = 1750
$ws.Range("M55").Value = -514.375
$ws.Range("N55").Value = -2096
$ws.Range("H93").Value = 2366.5
$ws.Range("I93").Value = 2366.3333
$ws.Range("K93").Value = 2366.3333
$ws.Range("M93").Value = -1118.3333
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 318
$ws.Range("H4").Value = 2020339
$ws.Range("J4").Value = 19173.75
$ws.Range("L4").Value = 19173.75
$ws.Range("N4").Value = -19399.75
$ws.Range("H64").Value = 65333.332
$ws.Range("J64").Value = 65333.332
$ws.Range("L64").Value = 65333.332
$ws.Range("N64").Value = -65829.33199999999
$ws.Range("H67").Value = 65333.332
$ws.Range("J67").Value = 65333.332
$ws.Range("L67").Value = 65333.332
$ws.Range("N67").Value = -67049.33199999999
$ws.Range("H113").Value = 1136.3
$ws.Range("I113").Value = 912.3333
$ws.Range("J113").Value = 1472.25
$ws.Range("K113").Value = 2736.9999
$ws.Range("L113").Value = 4416.75
$ws.Range("M113").Value = -566.9998999999998
$ws.Range("N113").Value = -8756.75
$ws.Range("H122").Value = 4337.25
$ws.Range("J122").Value = 2399.6667
$ws.Range("L122").Value = 7199.000100000001
$ws.Range("N122").Value = -12099.0001

